# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview sheet: status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets: the "Latest Target File", "Latest Handback File"
#    and "Latest Handback DateTime" columns (I/J/K) are now populated for
#    both data rows, with a new hyperlink on the target-file cell.
#  - A couple of columns are widened to fit the new content.

$wb = $excel.ActiveWorkbook

$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19486a4575ff0687249d96f5ce3ec19eef26ea4a/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-03 06:39:59"

$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-03 06:39:59"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")

$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(10).ColumnWidth = 39.1

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 06:40:13"

$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-03 06:40:13"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(10).ColumnWidth = 39.1
